# 17.1.1 - add two new reporting years (2021, 2022) as columns R and S,
# and refresh the figures for the last couple of existing years (P, Q)
# to match the newly-republished source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Create columns R and S by copying the formatting of column Q ---
# (column Q is the most recently added year column, so its formatting is
# the correct template for the two new year columns)
$ws.Range("Q4:Q10").Copy($ws.Range("R4:R10"))
$ws.Range("Q4:Q10").Copy($ws.Range("S4:S10"))

# --- 2. Header row: years ---
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# --- 3. Row 5: Revenues, total ---
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# --- 4. Row 6: Tax revenues ---
$ws.Range("P6").Value = 18.600000000000001
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# --- 5. Row 7: Received official transfers (no data -> "-") ---
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# --- 6. Row 8: Non-tax revenues ---
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# --- 7. Row 9: Revenues from the sale of non-financial assets ---
$ws.Range("P9").Value = 4.9000000000000004
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# --- 8. Row 10: Contributions / deductions for social needs ---
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# --- 9. Update the active selection to reflect where the author left off ---
$ws.Range("T3").Select() | Out-Null
